{"js": "// Bump the stack-trace line numbers embedded in the document text to\n// match the \"2.0.0 -> 2.0.1\" line-shift described in the commit.\n// Each pair is [searchText, replacementText]. Longer / more specific\n// strings are matched first implicitly because each search() call only\n// targets its own exact literal text, so there is no ambiguity between\n// e.g. \":278)\" and \":1278)\".\nconst replacements = [\n  [\"M2DocEvaluator.java:1120\", \"M2DocEvaluator.java:1132\"],\n  [\"M2DocEvaluator.java:1084\", \"M2DocEvaluator.java:1096\"],\n  [\"M2DocEvaluator.java:1300\", \"M2DocEvaluator.java:1305\"],\n  [\"M2DocEvaluator.java:278)\", \"M2DocEvaluator.java:283)\"],\n  [\"M2DocEvaluator.java:267)\", \"M2DocEvaluator.java:272)\"],\n  [\"AbstractTemplatesTestSuite.java:475\", \"AbstractTemplatesTestSuite.java:479\"],\n  [\"AbstractTemplatesTestSuite.java:384\", \"AbstractTemplatesTestSuite.java:388\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replacementText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Bump the stack-trace line numbers embedded in the document text to\n# match the \"2.0.0 -> 2.0.1\" line-shift described in the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"M2DocEvaluator.java:1120\", \"M2DocEvaluator.java:1132\"),\n    @(\"M2DocEvaluator.java:1084\", \"M2DocEvaluator.java:1096\"),\n    @(\"M2DocEvaluator.java:1300\", \"M2DocEvaluator.java:1305\"),\n    @(\"M2DocEvaluator.java:278)\", \"M2DocEvaluator.java:283)\"),\n    @(\"M2DocEvaluator.java:267)\", \"M2DocEvaluator.java:272)\"),\n    @(\"AbstractTemplatesTestSuite.java:475\", \"AbstractTemplatesTestSuite.java:479\"),\n    @(\"AbstractTemplatesTestSuite.java:384\", \"AbstractTemplatesTestSuite.java:388\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
